# Auto-generated update of the F column ('想去人数' / want-to-go count)
# values across the 展览, 演出 and 全部类型 sheets, reflecting a refreshed
# data pull (gh-pages output regeneration).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 180
$ws.Range("F3").Value = 102
$ws.Range("F5").Value = 953
$ws.Range("F6").Value = 5120
$ws.Range("F7").Value = 416
$ws.Range("F8").Value = 595
$ws.Range("F9").Value = 893
$ws.Range("F10").Value = 806
$ws.Range("F12").Value = 21
$ws.Range("F13").Value = 547
$ws.Range("F14").Value = 4
$ws.Range("F16").Value = 26
$ws.Range("F17").Value = 1687
$ws.Range("F18").Value = 1435
$ws.Range("F19").Value = 775
$ws.Range("F21").Value = 181
$ws.Range("F22").Value = 286
$ws.Range("F23").Value = 492
$ws.Range("F24").Value = 123
$ws.Range("F25").Value = 1042
$ws.Range("F28").Value = 2275
$ws.Range("F29").Value = 163
$ws.Range("F31").Value = 64
$ws.Range("F32").Value = 18
$ws.Range("F33").Value = 229
$ws.Range("F36").Value = 9
$ws.Range("F38").Value = 269
$ws.Range("F39").Value = 597
$ws.Range("F41").Value = 36
$ws.Range("F42").Value = 36

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 146
$ws.Range("F6").Value = 102

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 180
$ws.Range("F4").Value = 102
$ws.Range("F5").Value = 953
$ws.Range("F7").Value = 5120
$ws.Range("F8").Value = 416
$ws.Range("F9").Value = 595
$ws.Range("F11").Value = 146
$ws.Range("F12").Value = 893
$ws.Range("F13").Value = 806
$ws.Range("F15").Value = 102
$ws.Range("F17").Value = 21
$ws.Range("F18").Value = 547
$ws.Range("F19").Value = 4
$ws.Range("F22").Value = 26
$ws.Range("F23").Value = 1687
$ws.Range("F24").Value = 1435
$ws.Range("F25").Value = 775
$ws.Range("F27").Value = 181
$ws.Range("F28").Value = 286
$ws.Range("F30").Value = 492
$ws.Range("F31").Value = 123
$ws.Range("F32").Value = 1042
$ws.Range("F34").Value = 2276
$ws.Range("F35").Value = 163
$ws.Range("F37").Value = 18
$ws.Range("F38").Value = 229
$ws.Range("F41").Value = 9
$ws.Range("F42").Value = 269
$ws.Range("F43").Value = 597
$ws.Range("F45").Value = 36
$ws.Range("F46").Value = 36

Write-Host "Updated F-column view counts across sheets."
